# Apply corrected Qty/Amount (and a few rollup/subtotal and swapped-row)
# values to the "CryCompanywiseStockReport" stock sheet, per commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F24").Value = 146
$ws.Range("G24").Value = 16829.42
$ws.Range("F32").Value = 12
$ws.Range("G32").Value = 307.32
$ws.Range("F33").Value = 26
$ws.Range("G33").Value = 932.36
$ws.Range("F35").Value = 64
$ws.Range("G35").Value = 3278.72
$ws.Range("B40").Value = 53616.41
$ws.Range("F42").Value = 72
$ws.Range("G42").Value = 14167.44
$ws.Range("F47").Value = 192
$ws.Range("G47").Value = 37034.88
$ws.Range("F64").Value = 52
$ws.Range("G64").Value = 4052.36
$ws.Range("F67").Value = 196
$ws.Range("G67").Value = 51103.08
$ws.Range("B72").Value = 180010.83
$ws.Range("F84").Value = 92
$ws.Range("G84").Value = 2881.44
$ws.Range("B89").Value = 12037.43
$ws.Range("F111").Value = 68
$ws.Range("G111").Value = 5519.56
$ws.Range("F118").Value = 35
$ws.Range("G118").Value = 2764.3
$ws.Range("F120").Value = 18
$ws.Range("G120").Value = 841.3200000000001
$ws.Range("B129").Value = 68615.61
$ws.Range("F148").Value = 11
$ws.Range("G148").Value = 1092.41
$ws.Range("B153").Value = 19646.95
$ws.Range("F179").Value = 34
$ws.Range("G179").Value = 2672.74
$ws.Range("F180").Value = 37
$ws.Range("G180").Value = 5933.32
$ws.Range("F181").Value = 26
$ws.Range("G181").Value = 7540.26
$ws.Range("F183").Value = 270
$ws.Range("G183").Value = 12573.9
$ws.Range("F188").Value = 11
$ws.Range("G188").Value = 984.9400000000001
$ws.Range("B199").Value = 57995.52
$ws.Range("F212").Value = 50
$ws.Range("G212").Value = 3240
$ws.Range("B214").Value = 3240
$ws.Range("F221").Value = 154
$ws.Range("G221").Value = 17292.66
$ws.Range("B224").Value = 67929.42
$ws.Range("F228").Value = 418
$ws.Range("G228").Value = 7733
$ws.Range("B235").Value = 16420.59
$ws.Range("F238").Value = 22
$ws.Range("G238").Value = 2521.86
$ws.Range("B246").Value = 13146.08
$ws.Range("F261").Value = 6
$ws.Range("G261").Value = 1894.8
$ws.Range("F267").Value = 82
$ws.Range("G267").Value = 6724.82
$ws.Range("F272").Value = 6
$ws.Range("G272").Value = 399.54
$ws.Range("F273").Value = 108
$ws.Range("G273").Value = 4587.84
$ws.Range("F276").Value = 81
$ws.Range("G276").Value = 3943.89
$ws.Range("F277").Value = 17
$ws.Range("G277").Value = 1947.52
$ws.Range("F290").Value = 138
$ws.Range("G290").Value = 6468.06
$ws.Range("F291").Value = 3
$ws.Range("G291").Value = 332.82
$ws.Range("F295").Value = 13
$ws.Range("G295").Value = 1315.6
$ws.Range("B296").Value = 66194
$ws.Range("C296").Value = "HIM-Total Care Baby Pants Diapers-M-9s"
$ws.Range("F296").Value = 22
$ws.Range("G296").Value = 1884.96
$ws.Range("B297").Value = 64983
$ws.Range("C297").Value = "HIM-TOTAL CARE BABY PANTS DIAPERS-M-9S"
$ws.Range("F297").Value = 6
$ws.Range("G297").Value = 514.08
$ws.Range("B301").Value = 101942.25
$ws.Range("F312").Value = 25
$ws.Range("G312").Value = 3587
$ws.Range("F331").Value = 0
$ws.Range("G331").Value = 0
$ws.Range("F332").Value = 10
$ws.Range("G332").Value = 3205.2
$ws.Range("B334").Value = -22516.99
$ws.Range("F355").Value = 132
$ws.Range("G355").Value = 9851.16
$ws.Range("B362").Value = 74504.28999999999
$ws.Range("F372").Value = 45
$ws.Range("G372").Value = 2489.85
$ws.Range("F374").Value = 47
$ws.Range("G374").Value = 1504.94
$ws.Range("F376").Value = 179
$ws.Range("G376").Value = 29712.21
$ws.Range("B378").Value = 51521.07
$ws.Range("B387").Value = 58047
$ws.Range("D387").Value = 105.54
$ws.Range("E387").Value = 126.1
$ws.Range("F387").Value = 32
$ws.Range("G387").Value = 3377.28
$ws.Range("B388").Value = 47097
$ws.Range("D388").Value = 112.28
$ws.Range("E388").Value = 134.16
$ws.Range("F388").Value = 15
$ws.Range("G388").Value = 1684.2
$ws.Range("F393").Value = 377
$ws.Range("G393").Value = 36418.2
$ws.Range("B395").Value = 52567.74
$ws.Range("F425").Value = 20
$ws.Range("G425").Value = 2726.2
$ws.Range("F427").Value = 52
$ws.Range("G427").Value = 3989.96
$ws.Range("F429").Value = 51
$ws.Range("G429").Value = 4263.09
$ws.Range("B433").Value = 22563.21
$ws.Range("F436").Value = 205
$ws.Range("G436").Value = 9487.4
$ws.Range("F437").Value = 8
$ws.Range("G437").Value = 215.12
$ws.Range("B444").Value = 21829.05
$ws.Range("F456").Value = 6
$ws.Range("G456").Value = 1333.38
$ws.Range("F458").Value = 49
$ws.Range("G458").Value = 13289.78
$ws.Range("F459").Value = 26
$ws.Range("G459").Value = 3778.32
$ws.Range("B464").Value = 85561.2
$ws.Range("B485").Value = 53319
$ws.Range("E485").Value = 310.64
$ws.Range("F485").Value = -6
$ws.Range("G485").Value = -1643.52
$ws.Range("B486").Value = 64810
$ws.Range("E486").Value = 291.22
$ws.Range("F486").Value = 0
$ws.Range("G486").Value = 0
$ws.Range("B502").Value = 60025
$ws.Range("E502").Value = 37.22
$ws.Range("F502").Value = -98
$ws.Range("G502").Value = -3217.34
$ws.Range("B503").Value = 64833
$ws.Range("E503").Value = 34.9
$ws.Range("F503").Value = 88
$ws.Range("G503").Value = 2889.04
$ws.Range("F517").Value = 190
$ws.Range("G517").Value = 18975.3
$ws.Range("F528").Value = 53
$ws.Range("G528").Value = 1411.92
$ws.Range("B531").Value = 110899.2
$ws.Range("F537").Value = 183
$ws.Range("G537").Value = 6059.13
$ws.Range("F540").Value = 114
$ws.Range("G540").Value = 4988.64
$ws.Range("B541").Value = 20160.3
$ws.Range("F557").Value = 7
$ws.Range("G557").Value = 5216.33
$ws.Range("B562").Value = 38327.21
$ws.Range("F564").Value = 149
$ws.Range("G564").Value = 18155.65
$ws.Range("B567").Value = 20575.37
$ws.Range("F569").Value = 13
$ws.Range("G569").Value = 2429.18
$ws.Range("B579").Value = 14493.23
$ws.Range("F593").Value = 0
$ws.Range("G593").Value = 0
$ws.Range("B594").Value = 0
$ws.Range("F611").Value = 163
$ws.Range("G611").Value = 21695.3
$ws.Range("B613").Value = 21695.3
$ws.Range("F631").Value = 299
$ws.Range("G631").Value = 11012.17
$ws.Range("B634").Value = 195663
$ws.Range("F680").Value = 488
$ws.Range("G680").Value = 79597.67999999999
$ws.Range("B686").Value = 80610.23
$ws.Range("F696").Value = 22
$ws.Range("G696").Value = 3066.36
$ws.Range("B697").Value = 10276.19
$ws.Range("F710").Value = 6
$ws.Range("G710").Value = 2142.6
$ws.Range("F712").Value = 106
$ws.Range("G712").Value = 4154.14
$ws.Range("F713").Value = 128
$ws.Range("G713").Value = 4130.56
$ws.Range("B719").Value = 59366.38
$ws.Range("B724").Value = 2427168.36
$ws.Range("B725").Value = 2427168.36
